$wb = $excel.ActiveWorkbook

# Sheets "展览" and "全部类型" both contain the same data table; update the
# "想去人数" (want-to-go count) values in F2 and F3 on each.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 89
    $ws.Range("F3").Value = 13
}
